$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Hai-Xin%Bo%NULL%2,                  Wen%Li%NULL%1,                  Yuan%Yang%NULL%3,                  Yuan%Yang%NULL%0,                  Yu%Wang%NULL%0,                  Qinge%Zhang%NULL%1,                  Teris%Cheung%NULL%1,                  Xinjuan%Wu%NULL%1,                  Yu-Tao%Xiang%NULL%1]'
$ws.Range("I2").Value = '_PMC_CROSSREF'
$ws.Range("E3").Value = '[Haozheng%Cai%NULL%2,                  Baoren%Tu%NULL%2,                  Jing%Ma%NULL%2,                  Limin%Chen%NULL%2,                  Lei%Fu%NULL%2,                  Yongfang%Jiang%NULL%2,                  Quan%Zhuang%NULL%2]'
$ws.Range("I3").Value = '_PMC_CROSSREF'
$ws.Range("E4").Value = '[Jinya%Cao%NULL%1,                  Jing%Wei%NULL%1,                  Huadong%Zhu%NULL%1,                  Yanping%Duan%NULL%1,                  Wenqi%Geng%NULL%1,                  Xia%Hong%NULL%1,                  Jing%Jiang%NULL%1,                  Xiaohui%Zhao%NULL%1,                  Boheng%Zhu%NULL%1]'
$ws.Range("I4").Value = '_PMC_CROSSREF'
$ws.Range("E5").Value = '[Wenjun%Cao%wjcao16@czmc.edu.cn%2,                  Ziwei%Fang%NULL%2,                  Guoqiang%Hou%NULL%2,                  Mei%Han%NULL%2,                  Xinrong%Xu%NULL%2,                  Jiaxin%Dong%NULL%2,                  Jianzhong%Zheng%zhengjj@czmc.edu.cn%2]'
$ws.Range("I5").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E6").Value = '[Joseph PY%Chung%xref no email%2,            Wai-song%Yeung%xref no email%2]'
$ws.Range("I6").Value = '_CROSSREF'
$ws.Range("E7").Value = '[Joseph PY%Chung%xref no email%0,            Wai-song%Yeung%xref no email%0]'
$ws.Range("I7").Value = '_CROSSREF'
$ws.Range("E8").Value = '[Junling%Gao%NULL%1,                  Pinpin%Zheng%NULL%2,                  Pinpin%Zheng%NULL%0,                  Yingnan%Jia%NULL%1,                  Hao%Chen%NULL%1,                  Yimeng%Mao%NULL%1,                  Suhong%Chen%NULL%1,                  Yi%Wang%NULL%1,                  Hua%Fu%NULL%1,                  Junming%Dai%NULL%1,                  Kenji%Hashimoto%NULL%10,                  Kenji%Hashimoto%NULL%0,                  Junling%Gao%NULL%3,                  Junling%Gao%NULL%0,                  Junling%Gao%NULL%0,                  Kenji%Hashimoto%NULL%0,                  Kenji%Hashimoto%NULL%0,                  Kenji%Hashimoto%NULL%0,                  Kenji%Hashimoto%NULL%0,                  Kenji%Hashimoto%NULL%0]'
$ws.Range("I8").Value = '_PMC_CROSSREF'
$ws.Range("E10").Value = '[ Yueqin%Huang%null%1,                   Yu%Wang%null%0,                   Hong%Wang%null%1,                   Zhaorui%Liu%null%1,                   Xin%Yu%null%1,                   Jie%Yan%null%2,                   Yaqin%Yu%null%1,                   Changgui%Kou%null%1,                   Xiufeng%Xu%null%1,                   Jin%Lu%null%1,                   Zhizhong%Wang%null%1,                   Shulan%He%null%1,                   Yifeng%Xu%null%1,                   Yanling%He%null%1,                   Tao%Li%null%3,                   Wanjun%Guo%null%1,                   Hongjun%Tian%null%1,                   Guangming%Xu%null%1,                   Xiangdong%Xu%null%1,                   Yanjuan%Ma%null%1,                   Linhong%Wang%null%1,                   Limin%Wang%null%1,                   Yongping%Yan%null%1,                   Bo%Wang%null%0,                   Shuiyuan%Xiao%null%1,                   Liang%Zhou%null%1,                   Lingjiang%Li%null%1,                   Liwen%Tan%null%1,                   Tingting%Zhang%null%1,                   Chao%Ma%null%1,                   Qiang%Li%null%2,                   Hua%Ding%null%1,                   Hongchun%Geng%null%1,                   Fujun%Jia%null%1,                   Jianfei%Shi%null%1,                   Shiliang%Wang%null%1,                   Ning%Zhang%null%1,                   Xinbai%Du%null%1,                   Xiangdong%Du%null%1,                   Yue%Wu%null%2,                  Yueqin%Huang%null%1,                  Yu%Wang%null%0,                  Hong%Wang%null%1,                  Zhaorui%Liu%null%1,                  Xin%Yu%null%1,                  Jie%Yan%null%1,                  Yaqin%Yu%null%1,                  Changgui%Kou%null%1,                  Xiufeng%Xu%null%1,                  Jin%Lu%null%1,                  Zhizhong%Wang%null%1,                  Shulan%He%null%1,                  Yifeng%Xu%null%2,                  Yanling%He%null%1,                  Tao%Li%null%0,                  Wanjun%Guo%null%1,                  Hongjun%Tian%null%1,                  Guangming%Xu%null%1,                  Xiangdong%Xu%null%1,                  Yanjuan%Ma%null%1,                  Linhong%Wang%null%1,                  Limin%Wang%null%1,                  Yongping%Yan%null%1,                  Bo%Wang%null%1,                  Shuiyuan%Xiao%null%1,                  Liang%Zhou%null%1,                  Lingjiang%Li%null%1,                  Liwen%Tan%null%1,                  Tingting%Zhang%null%1,                  Chao%Ma%null%1,                  Qiang%Li%null%1,                  Hua%Ding%null%1,                  Hongchun%Geng%null%1,                  Fujun%Jia%null%1,                  Jianfei%Shi%null%1,                  Shiliang%Wang%null%1,                  Ning%Zhang%null%0,                  Xinbai%Du%null%1,                  Xiangdong%Du%null%1,                  Yue%Wu%null%1]'
$ws.Range("I10").Value = '_CROSSREF'
$ws.Range("E11").Value = '[Lijun%Kang%NULL%2,                  Simeng%Ma%NULL%4,                  Min%Chen%NULL%2,                  Jun%Yang%NULL%2,                  Ying%Wang%NULL%4,                  Ruiting%Li%NULL%5,                  Lihua%Yao%NULL%4,                  Hanping%Bai%NULL%2,                  Zhongxiang%Cai%NULL%4,                  Bing%Xiang Yang%NULL%2,                  Shaohua%Hu%NULL%4,                  Kerang%Zhang%NULL%2,                  Gaohua%Wang%NULL%4,                  Ci%Ma%NULL%2,                  Zhongchun%Liu%NULL%4]'
$ws.Range("I11").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E12").Value = '[Kin On%Kwok%NULL%0,                  Kin Kit%Li%NULL%1,                  Henry Ho Hin%Chan%NULL%1,                  Yuan Yuan%Yi%NULL%1,                  Arthur%Tang%NULL%1,                  Wan In%Wei%NULL%1,                  Samuel Yeung Shan%Wong%NULL%1]'
$ws.Range("I12").Value = '_PMC_CROSSREF'
$ws.Range("E13").Value = '[Jianbo%Lai%NULL%5,                  Simeng%Ma%NULL%0,                  Ying%Wang%NULL%0,                  Zhongxiang%Cai%NULL%0,                  Jianbo%Hu%NULL%2,                  Ning%Wei%NULL%2,                  Jiang%Wu%NULL%2,                  Hui%Du%NULL%2,                  Tingting%Chen%NULL%2,                  Ruiting%Li%NULL%0,                  Huawei%Tan%NULL%2,                  Lijun%Kang%NULL%2,                  Lihua%Yao%NULL%0,                  Manli%Huang%NULL%2,                  Huafen%Wang%NULL%2,                  Gaohua%Wang%NULL%0,                  Zhongchun%Liu%NULL%0,                  Shaohua%Hu%NULL%0]'
$ws.Range("I13").Value = '_PMC_CROSSREF'
$ws.Range("E14").Value = '[Nianqi%Liu%NULL%2,                  Fan%Zhang%NULL%1,                  Cun%Wei%NULL%1,                  Yanpu%Jia%NULL%1,                  Zhilei%Shang%NULL%1,                  Luna%Sun%NULL%1,                  Lili%Wu%NULL%1,                  Zhuoer%Sun%NULL%1,                  Yaoguang%Zhou%NULL%1,                  Yan%Wang%NULL%1,                  Weizhi%Liu%13024141970@163.com%1]'
$ws.Range("I14").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E15").Value = '[Wen%Lu%luwen67@sina.com%3,                  Hang%Wang%NULL%1,                  Yuxing%Lin%NULL%1,                  Li%Li%lilifuzhou@126.com%0]'
$ws.Range("I15").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E16").Value = '[Yuanyuan%Mo%NULL%2,                  Lan%Deng%NULL%4,                  Lan%Deng%NULL%0,                  Liyan%Zhang%NULL%2,                  Qiuyan%Lang%NULL%2,                  Chunyan%Liao%NULL%2,                  Nannan%Wang%NULL%2,                  Mingqin%Qin%2026142822@qq.com%2,                  Huiqiao%Huang%820325832@qq.com%2]'
$ws.Range("I16").Value = '_PMC_CROSSREF'
$ws.Range("I17").Value = '_CROSSREF'
$ws.Range("E18").Value = '[Hoang C.%Nguyen%NULL%1,                  Minh H.%Nguyen%NULL%1,                  Binh N.%Do%NULL%1,                  Cuong Q.%Tran%NULL%2,                  Cuong Q.%Tran%NULL%0,                  Thao T. P.%Nguyen%NULL%1,                  Khue M.%Pham%NULL%2,                  Khue M.%Pham%NULL%0,                  Linh V.%Pham%NULL%2,                  Linh V.%Pham%NULL%0,                  Khanh V.%Tran%NULL%1,                  Trang T.%Duong%NULL%1,                  Tien V.%Tran%NULL%1,                  Thai H.%Duong%NULL%1,                  Tham T.%Nguyen%NULL%1,                  Quyen H.%Nguyen%NULL%1,                  Thanh M.%Hoang%NULL%1,                  Kien T.%Nguyen%NULL%1,                  Thu T. M.%Pham%NULL%1,                  Shwu-Huey%Yang%NULL%1,                  Jane C.-J.%Chao%NULL%2,                  Jane C.-J.%Chao%NULL%0,                  Tuyen Van%Duong%NULL%1]'
$ws.Range("I18").Value = '_PMC_CROSSREF'
$ws.Range("E19").Value = '[Federica%Pulvirenti%NULL%1,                  Francesco%Cinetto%NULL%1,                  Cinzia%Milito%NULL%1,                  Livia%Bonanni%NULL%1,                  Anna Maria%Pesce%NULL%1,                  Giorgia%Leodori%NULL%1,                  Giulia%Garzi%NULL%1,                  Marzia%Miglionico%NULL%1,                  Stefano%Tabolli%NULL%1,                  Isabella%Quinti%isabella.quinti@uniroma1.it%1]'
$ws.Range("I19").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E20").Value = '[Jianyin%Qiu%NULL%1,                  Bin%Shen%NULL%1,                  Min%Zhao%NULL%1,                  Zhen%Wang%NULL%2,                  Zhen%Wang%NULL%0,                  Bin%Xie%NULL%1,                  Yifeng%Xu%NULL%0]'
$ws.Range("I20").Value = '_PMC_CROSSREF'
$ws.Range("E21").Value = '[Deblina%Roy%NULL%1,                  Sarvodaya%Tripathy%NULL%1,                  Sujita Kumar%Kar%NULL%1,                  Nivedita%Sharma%NULL%1,                  Sudhir Kumar%Verma%NULL%1,                  Vikas%Kaushal%NULL%1]'
$ws.Range("I21").Value = '_PMC_elsevier_CROSSREF'
$ws.Range("E22").Value = '[Benjamin Y.Q.%Tan%NULL%2,                  Nicholas W.S.%Chew%NULL%1,                  Grace K.H.%Lee%NULL%1,                  Mingxue%Jing%NULL%1,                  Yihui%Goh%NULL%1,                  Leonard L.L.%Yeo%NULL%1,                  Ka%Zhang%NULL%1,                  Howe-Keat%Chin%NULL%1,                  Aftab%Ahmad%NULL%1,                  Faheem Ahmed%Khan%NULL%1,                  Ganesh Napolean%Shanmugam%NULL%1,                  Bernard P.L.%Chan%NULL%1,                  Sibi%Sunny%NULL%1,                  Bharatendu%Chandra%NULL%1,                  Jonathan J.Y.%Ong%NULL%1,                  Prakash R.%Paliwal%NULL%1,                  Lily Y.H.%Wong%NULL%1,                  Renarebecca%Sagayanathan%NULL%1,                  Jin Tao%Chen%NULL%1,                  Alison Ying%Ying Ng%NULL%1,                  Hock Luen%Teoh%NULL%1,                  Cyrus S.%Ho%NULL%3,                  Roger C.%Ho%NULL%5,                  Vijay K.%Sharma%NULL%2]'
$ws.Range("I22").Value = '_PMC_CROSSREF'
$ws.Range("E23").Value = '[Cuiyan%Wang%NULL%4,                  Riyu%Pan%NULL%5,                  Riyu%Pan%NULL%0,                  Xiaoyang%Wan%NULL%3,                  Yilin%Tan%NULL%3,                  Linkang%Xu%NULL%3,                  Cyrus S.%Ho%NULL%0,                  Roger C.%Ho%NULL%0,                  Roger C.%Ho%NULL%0]'
$ws.Range("I23").Value = '_PMC_CROSSREF'
$ws.Range("E24").Value = '[Yenan%Wang%xref no email%2,            Yu%Di%xref no email%2,            Junjie%Ye%xref no email%2,            Wenbin%Wei%xref no email%2]'
$ws.Range("I24").Value = '_CROSSREF'
$ws.Range("E25").Value = '[Michael S.%Wolf%NULL%1,                  Marina%Serper%NULL%1,                  Lauren%Opsasnick%NULL%1,                  Rachel M.%O''Conor%NULL%1,                  Laura M.%Curtis%NULL%1,                  Julia Yoshino%Benavente%NULL%1,                  Guisselle%Wismer%NULL%1,                  Stephanie%Batio%NULL%1,                  Morgan%Eifler%NULL%1,                  Pauline%Zheng%NULL%1,                  Andrea%Russell%NULL%1,                  Marina%Arvanitis%NULL%1,                  Daniela%Ladner%NULL%1,                  Mary%Kwasny%NULL%1,                  Stephen D.%Persell%NULL%1,                  Theresa%Rowe%NULL%1,                  Jeffrey A.%Linder%NULL%1,                  Stacy C.%Bailey%NULL%1]'
$ws.Range("I25").Value = '_PMC_CROSSREF'
$ws.Range("E26").Value = '[Han%Xiao%NULL%2,                  Yan%Zhang%NULL%2,                  Desheng%Kong%NULL%2,                  Shiyue%Li%NULL%2,                  Ningxi%Yang%NULL%2]'
$ws.Range("I26").Value = '_PMC_CROSSREF'
$ws.Range("E27").Value = '[Shuai%Yuan%NULL%2,                  Zhenxin%Liao%NULL%2,                  Haojie%Huang%NULL%2,                  Boyue%Jiang%NULL%2,                  Xueyan%Zhang%NULL%2,                  Yingwen%Wang%NULL%2,                  Mingyi%Zhao%NULL%2]'
$ws.Range("I27").Value = '_PMC_CROSSREF'
$ws.Range("E28").Value = '[Wen-rui%Zhang%NULL%3,                  Kun%Wang%NULL%0,                  Lu%Yin%NULL%0,                  Wen-feng%Zhao%NULL%1,                  Qing%Xue%NULL%1,                  Mao%Peng%NULL%1,                  Bao-quan%Min%NULL%1,                  Qing%Tian%NULL%1,                  Hai-xia%Leng%NULL%1,                  Jia-lin%Du%NULL%1,                  Hong%Chang%NULL%1,                  Yuan%Yang%NULL%0,                  Wei%Li%NULL%1,                  Fang-fang%Shangguan%NULL%1,                  Tian-yi%Yan%NULL%1,                  Hui-qing%Dong%NULL%1,                  Ying%Han%NULL%1,                  Yu-ping%Wang%NULL%1,                  Fiammetta%Cosci%NULL%1,                  Hong-xing%Wang%NULL%1]'
$ws.Range("I28").Value = '_PMC_CROSSREF'
$ws.Range("E29").Value = '[Yingfei%Zhang%NULL%4,                  Zheng Feei%Ma%NULL%4,                  Zheng Feei%Ma%NULL%0]'
$ws.Range("I29").Value = '_PMC_CROSSREF'
$ws.Range("E30").Value = '[M.%Zhang%NULL%1,                  M.%Zhou%NULL%1,                  F.%Tang%NULL%1,                  Y.%Wang%NULL%1,                  H.%Nie%NULL%1,                  L.%Zhang%NULL%1,                  G.%You%NULL%1]'
$ws.Range("I30").Value = '_PMC_elsevier_CROSSREF'
